$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update odds on row 5 (Operario vs Sport Recife) ---
$ws.Range("G5").Value = 3.6
$ws.Range("I5").Value = 2.15
$ws.Range("J5").Value = 4.5
$ws.Range("L5").Value = 3
$ws.Range("X5").Value = 17
$ws.Range("AB5").Value = 51
$ws.Range("AF5").Value = 81
$ws.Range("AJ5").Value = 19

# --- Update odds on row 6 (Ponte Preta vs Paysandu PA) ---
$ws.Range("G6").Value = 2.45
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 2.9
$ws.Range("J6").Value = 3.25
$ws.Range("Q6").Value = 2.35
$ws.Range("R6").Value = 1.57
$ws.Range("AA6").Value = 23
$ws.Range("AG6").Value = 7.5
$ws.Range("AJ6").Value = 29
$ws.Range("AP6").Value = 29
$ws.Range("AS6").Value = 251

# --- Update odds on row 7 (Ituano vs CRB) ---
$ws.Range("AM7").Value = 1000

# --- Update odds on row 8 (Atl. Nacional vs Santa Fe) ---
$ws.Range("V8").Value = 1.54

# --- Remove the match that no longer belongs in the sheet ---
# Row 9: Sol de America vs Nacional Asuncion (PARAGUAY - PRIMERA DIVISION)
# Deleting the entire row shifts rows 10 and 11 up to 9 and 10.
$ws.Rows.Item(9).Delete()
